# Append three new geocoded address/coordinate rows to the slave data sheet,
# continuing directly after the existing 4 rows (rows 5, 6 and 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Frederiksborgvej 13-5  "
$ws.Range("B5").Value = "12.531159, 55.71852"

$ws.Range("A6").Value = "Buddingevej 66  "
$ws.Range("B6").Value = "12.495893, 55.762853"

$ws.Range("A7").Value = "Ny Kronborgvej 2  "
$ws.Range("B7").Value = "12.614004, 56.039332"
